# Auto-generated edit script: update cryptos list per commit
# "Updated cryptos list on Fri Mar 10 07:07:08 UTC 2023 with GitHub Actions"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "19.919.46"
$ws.Range("E2").Value = "  -8.34%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.410.01"
$ws.Range("E3").Value = "  -8.47%  "

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.004"
$ws.Range("E4").Value = "  +0.27%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "1.003"
$ws.Range("E5").Value = "  +0.28%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "272.57"
$ws.Range("E6").Value = "  -6.04%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3700"
$ws.Range("E7").Value = "  -4.66%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3061"
$ws.Range("E8").Value = "  -3.98%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "39.04"
$ws.Range("E9").Value = "  -9.51%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.9928"
$ws.Range("E10").Value = "  -6.49%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.06545"
$ws.Range("E11").Value = "  -9.16%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.004"
$ws.Range("E12").Value = "  +0.24%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.409"
$ws.Range("E13").Value = "  -4.38%  "

# Row 14
$ws.Range("B14").Value = "Chainlink"
$ws.Range("C14").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.150"
$ws.Range("E14").Value = "  -7.30%  "

# Row 15
$ws.Range("B15").Value = "Solana"
$ws.Range("C15").Value = "https://coinranking.com/coin/zNZHO_Sjf+solana-sol"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "16.78"
$ws.Range("E15").Value = "  -10.22%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.415.84"
$ws.Range("E16").Value = "  -8.06%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.00001004"
$ws.Range("E17").Value = "  -9.38%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.05732"
$ws.Range("E18").Value = "  -13.05%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "73.70"
$ws.Range("E19").Value = "  -11.50%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "1.003"
$ws.Range("E20").Value = "  +0.28%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "5.565"
$ws.Range("E21").Value = "  -9.73%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "14.36"
$ws.Range("E22").Value = "  -6.83%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "10.78"
$ws.Range("E23").Value = "  -1.22%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.325"
$ws.Range("E24").Value = "  -3.42%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "19.933.86"
$ws.Range("E25").Value = "  -8.32%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.250"
$ws.Range("E26").Value = "  -5.35%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "138.89"
$ws.Range("E27").Value = "  -5.23%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "16.87"
$ws.Range("E28").Value = "  -8.46%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.574.66"
$ws.Range("E29").Value = "  -8.14%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "109.01"
$ws.Range("E30").Value = "  -7.40%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.851"
$ws.Range("E31").Value = "  -20.37%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "5.340"
$ws.Range("E32").Value = "  -10.04%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.8517"
$ws.Range("E33").Value = "  -12.76%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.07684"
$ws.Range("E34").Value = "  -6.35%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "8.383"
$ws.Range("E35").Value = "  -5.22%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.05732"
$ws.Range("E36").Value = "  -6.09%  "

# Row 37
$ws.Range("E37").Value = "  +0.21%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "4.765"
$ws.Range("E38").Value = "  -7.38%  "

# Row 39
$ws.Range("B39").Value = "Algorand"
$ws.Range("C39").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.1909"
$ws.Range("E39").Value = "  -6.56%  "

# Row 40
$ws.Range("B40").Value = "VeChain"
$ws.Range("C40").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.02026"
$ws.Range("E40").Value = "  -8.13%  "

# Row 41
$ws.Range("B41").Value = "Aptos"
$ws.Range("C41").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "10.35"
$ws.Range("E41").Value = "  -3.23%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.053"
$ws.Range("E42").Value = "  -11.63%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.280"
$ws.Range("E43").Value = "  -13.09%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.5266"
$ws.Range("E44").Value = "  -8.51%  "

# Row 45
$ws.Range("B45").Value = "PancakeSwap"
$ws.Range("C45").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "3.522"
$ws.Range("E45").Value = "  -5.86%  "

# Row 46
$ws.Range("B46").Value = "EnergySwap"
$ws.Range("C46").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "12.21"
$ws.Range("E46").Value = "  -7.18%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.5087"
$ws.Range("E47").Value = "  -7.75%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.789"
$ws.Range("E48").Value = "  -4.74%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "109.05"
$ws.Range("E49").Value = "  -6.72%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.044"
$ws.Range("E50").Value = "  -10.16%  "

# Row 51
$ws.Range("E51").Value = "  +0.35%  "
